# Update the build version / timestamp strings that appear throughout the
# workbook after a new release build.
#
# Old build timestamp: "January 30 2026 16.19.47 EST"
# New build timestamp: "February 02 2026 12.49.33 EST"

$wb = $excel.ActiveWorkbook

$oldStamp = "January 30 2026 16.19.47 EST"
$newStamp = "February 02 2026 12.49.33 EST"

# --- "About" sheet -------------------------------------------------------
$aboutWs = $wb.Worksheets.Item("About")

$a2 = $aboutWs.Range("A2").Value()
$a2 = $a2.Replace($oldStamp, $newStamp)
$aboutWs.Range("A2").Value = $a2

$a6 = $aboutWs.Range("A6").Value()
$a6 = $a6.Replace($oldStamp, $newStamp)
$aboutWs.Range("A6").Value = $a6

# --- "Boundaries and methane sources" sheet ------------------------------
$dataWs = $wb.Worksheets.Item("Boundaries and methane sources")

# Column S holds the "build_version" value for each data row (rows 2-10).
for ($row = 2; $row -le 10; $row++) {
    $cell = $dataWs.Cells.Item($row, 19)  # column S = 19
    $val = $cell.Value()
    $val = $val.Replace($oldStamp, $newStamp)
    $cell.Value = $val
}
